# Refresh the cryptos price/volume table (GitHub Actions scheduled update).
# Column D ("Price") values are stored as literal text in the workbook, so a
# leading apostrophe is used to force Excel to keep them as text instead of
# auto-converting to numbers (which would drop significant trailing zeros,
# e.g. "0.110" -> 0.11, or mangle multi-dot values like "2.516.30").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.517.10"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "'2.514.65"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'317.87"
$ws.Range("E5").Value = "  +5.52%  "
$ws.Range("D6").Value = "'94.24"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'0.528"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").Value = "'35.63"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "'0.0813"
$ws.Range("E11").Value = "  +1.72%  "
$ws.Range("D12").Value = "'7.52"
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("E13").Value = "  -2.16%  "
$ws.Range("D14").Value = "'2.902.06"
$ws.Range("E14").Value = "  +0.79%  "
$ws.Range("D15").Value = "'2.552.17"
$ws.Range("E15").Value = "  +2.44%  "
$ws.Range("E16").Value = "  +2.81%  "
$ws.Range("D17").Value = "'0.846"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("D18").Value = "'42.622.99"
$ws.Range("E18").Value = "  +1.44%  "
$ws.Range("D19").Value = "'12.89"
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").Value = "'6.65"
$ws.Range("E20").Value = "  +5.00%  "
$ws.Range("D21").Value = "'0.0₃0957"
$ws.Range("E21").Value = "  -0.06%  "
$ws.Range("D22").Value = "'69.27"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").Value = "'250.66"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").Value = "'2.95"
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("D26").Value = "'26.67"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("E28").Value = "  +6.35%  "
$ws.Range("D29").Value = "'40.91"
$ws.Range("E29").Value = "  +11.10%  "
$ws.Range("D30").Value = "'10.22"
$ws.Range("E30").Value = "  +2.10%  "
$ws.Range("D31").Value = "'5.91"
$ws.Range("E31").Value = "  +1.47%  "
$ws.Range("D32").Value = "'156.98"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'2.11"
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("B34").Value = "Celestia"
$ws.Range("C34").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D34").Value = "'19.06"
$ws.Range("E34").Value = "  +5.31%  "
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("E36").Value = "  +0.57%  "
$ws.Range("D37").Value = "'0.0778"
$ws.Range("E37").Value = "  +0.73%  "
$ws.Range("D38").Value = "'0.110"
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").Value = "'23.39"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("D41").Value = "'2.28"
$ws.Range("E41").Value = "  +15.65%  "
$ws.Range("E42").Value = "  +0.12%  "
$ws.Range("D43").Value = "'0.0302"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("B44").Value = "NEARProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D44").Value = "'3.30"
$ws.Range("E44").Value = "  -0.30%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'3.76"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").Value = "'2.012.06"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'84.90"
$ws.Range("E47").Value = "  +2.94%  "
$ws.Range("D48").Value = "'8.92"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("D49").Value = "'74.51"
$ws.Range("E49").Value = "  +5.05%  "
$ws.Range("D50").Value = "'2.756.97"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "'101.92"
$ws.Range("E51").Value = "  +2.39%  "

Write-Output "Applied 96 cell updates"
